$wb = $excel.ActiveWorkbook

# zh-cn sheet: row 5 "Correspond Handoff Datetime" (D5) and "Correspond Handback DateTime" (G5)
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("D5").Value = "2016-01-18 07:04:54"
$wsZh.Range("G5").Value = "2016-01-18 07:05:38"

# de-de sheet: row 5 "Correspond Handoff Datetime" (D5) and "Correspond Handback DateTime" (G5)
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("D5").Value = "2016-01-18 07:05:04"
$wsDe.Range("G5").Value = "2016-01-18 07:05:54"
